# Add two new header columns (G, H) to Sheet1 and populate the new
# "weekly mean" / "Tx # SF/mL new" figures that go with them, per the
# commit "improving a few figures and working on paper".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column headers (row 1)
$ws.Range("G1").Value = "total crypto SF/mL weekly mean"
$ws.Range("H1").Value = "Tx # SF/mL new"

# Row 4
#   E4 gets re-entered (formula/value unchanged) so it drops back to the
#   sheet's default style, matching the rest of the table.
$ws.Range("E4").Clear()
$ws.Range("E4").Formula = "=(75166.434)/1000"
$ws.Range("G4").Formula = "=(520876.1/1000)"
$ws.Range("H4").Formula = "=(0.1788893*G4)"

# Row 8
$ws.Range("G8").Formula = "=(81760.8/1000)"
$ws.Range("H8").Formula = "=(D8*G8)"

# Row 10
$ws.Range("E10").Clear()
$ws.Range("E10").Formula = "=(352786.84)/1000"
$ws.Range("G10").Formula = "=(94497.04/1000)"
$ws.Range("H10").Formula = "=D10*G10"

# Row 15
$ws.Range("G15").Formula = "=(236434.8/1000)"
$ws.Range("H15").Formula = "=D15*G15"

# Leave the cursor where the author last worked, like the saved file.
$ws.Range("H16").Select() | Out-Null
